$p = $ppt.ActivePresentation

# --- Slide 2: "Steps, part 1" -------------------------------------------------
# Created by duplicating slide 1 (preserves shape ids / creationIds exactly
# like PowerPoint's "Duplicate Slide"), then editing the title & body text
# and repositioning the content placeholder.
$dup1 = $p.Slides.Item(1).Duplicate()
$slide2 = $dup1.Item(1)

$slide2.Shapes.Item(1).TextFrame.TextRange.Text = "Steps, part 1"

$body2 = $slide2.Shapes.Item(2)
$body2.Left = 66.0
$body2.Top = 119.22582677165354
$body2.Width = 828.0
$body2.Height = 367.1492156982422

$tr2 = $body2.TextFrame.TextRange
$tr2.Text = "Figure out a shader that only draw non-lit object`r" + `
    "Could be wire-frame or not`r" + `
    "Add a class that holds a list of debug objects to draw this frame`r" + `
    "Each object can have: `r" + `
    "Draw for 1 frame only `r" + `
    "Draw for ever`r" + `
    "Draw for a certain amount of time, then disappear`r" + `
    "Load these objects into their own VAO + shader combo`r" + `
    "Each frame, we go through this list of debug objects and draw then, then delete them if they aren’t needed anymore."

$tr2.Paragraphs(2,1).IndentLevel = 2
$tr2.Paragraphs(4,1).IndentLevel = 2
$tr2.Paragraphs(5,1).IndentLevel = 3
$tr2.Paragraphs(6,1).IndentLevel = 3
$tr2.Paragraphs(7,1).IndentLevel = 3

# --- Slide 3: "Steps, part 2 : any dynamic lines" -----------------------------
$dup2 = $p.Slides.Item(1).Duplicate()
$slide3 = $dup2.Item(1)
$slide3.MoveTo(3)

$slide3.Shapes.Item(1).TextFrame.TextRange.Text = "Steps, part 2 : any dynamic lines"

$body3 = $slide3.Shapes.Item(2)
$body3.Left = 66.0
$body3.Top = 119.22582677165354
$body3.Width = 828.0
$body3.Height = 367.1492156982422

$tr3 = $body3.TextFrame.TextRange
$tr3.Text = "Make a large vertex buffer (NOT an element buffer) that’s large enough to draw the maximum number of lines.`r" + `
    "Keep the original C++/CPU side array.`r" + `
    "Load the lines into this C++/CPU side array.`r" + `
    "Each frame, update the buffer on the GPU side.`r" + `
    "Tell it to draw X number of lines that are in that buffer"
